$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.135.09"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "'3.168.68"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'603.66"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'154.02"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'3.165.84"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").Value = "'0.546"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").Value = "'5.59"
$ws.Range("E11").Value = "  -9.33%  "
$ws.Range("D12").Value = "'0.518"
$ws.Range("E12").Value = "  +2.56%  "
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'38.37"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "'3.684.18"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "'66.156.19"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'7.44"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "'3.168.10"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").Value = "'510.92"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "'15.39"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'14.73"
$ws.Range("E24").Value = "  -2.99%  "
$ws.Range("D25").Value = "'84.63"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "'9.20"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").Value = "'2.39"
$ws.Range("E29").Value = "  +7.07%  "
$ws.Range("E30").Value = "  +7.92%  "
$ws.Range("D31").Value = "'7.16"
$ws.Range("E31").Value = "  +6.12%  "
$ws.Range("D32").Value = "'28.01"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").Value = "'6.50"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "'502.81"
$ws.Range("E36").Value = "  +5.25%  "
$ws.Range("D37").Value = "'54.74"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("D38").Value = "'0.0884"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").Value = "'0.0420"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'0.129"
$ws.Range("E40").Value = "  +8.67%  "
$ws.Range("D41").Value = "'8.77"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").Value = "'0.0₃0681"
$ws.Range("E42").Value = "  +7.05%  "
$ws.Range("D43").Value = "'0.296"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("E44").Value = "  -4.98%  "
$ws.Range("D45").Value = "'2.43"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'2.821.61"
$ws.Range("E46").Value = "  -3.72%  "
$ws.Range("D47").Value = "'27.92"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").Value = "'2.38"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").Value = "'35.20"
$ws.Range("E51").Value = "  +3.75%  "
